$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'244.93"
$c.ClearFormats()
$c = $ws.Range("D3")
$c.Value = "'23.86"
$c.ClearFormats()
$c = $ws.Range("D4")
$c.Value = "'5.179"
$c.ClearFormats()
$c = $ws.Range("D5")
$c.Value = "'0.05718"
$c.ClearFormats()
$c = $ws.Range("D6")
$c.Value = "'6.474"
$c.ClearFormats()
$c = $ws.Range("D7")
$c.Value = "'3.166"
$c.ClearFormats()
$c = $ws.Range("D8")
$c.Value = "'0.8128"
$c.ClearFormats()
$c = $ws.Range("D9")
$c.Value = "'0.8548"
$c.ClearFormats()
$c = $ws.Range("D10")
$c.Value = "'0.1371"
$c.ClearFormats()
$c = $ws.Range("D11")
$c.Value = "'0.06935"
$c.ClearFormats()
$c = $ws.Range("D12")
$c.Value = "'0.03190"
$c.ClearFormats()
$c = $ws.Range("D13")
$c.Value = "'0.02897"
$c.ClearFormats()
$c = $ws.Range("D14")
$c.Value = "'0.09331"
$c.ClearFormats()
$c = $ws.Range("D15")
$c.Value = "'3.812"
$c.ClearFormats()
$c = $ws.Range("D16")
$c.Value = "'0.001534"
$c.ClearFormats()
$c = $ws.Range("D17")
$c.Value = "'0.04706"
$c.ClearFormats()
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$c = $ws.Range("D18")
$c.Value = "'0.0005989"
$c.ClearFormats()
$ws.Range("E18").Value = "17OneONE"
$ws.Range("B19").Value = "TigerCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$c = $ws.Range("D19")
$c.Value = "'0.006223"
$c.ClearFormats()
$ws.Range("E19").Value = "18TigerCashTCH"
$ws.Range("B20").Value = "BitKan"
$ws.Range("C20").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$c = $ws.Range("D20")
$c.Value = "'0.001241"
$c.ClearFormats()
$ws.Range("E20").Value = "19BitKanKAN"
$ws.Range("B21").Value = "HotbitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$c = $ws.Range("D21")
$c.Value = "'0.004809"
$c.ClearFormats()
$ws.Range("E21").Value = "20HotbitTokenHTB"
$ws.Range("B22").Value = "NitroEx"
$ws.Range("C22").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$c = $ws.Range("D22")
$c.Value = "'0.00008495"
$c.ClearFormats()
$ws.Range("E22").Value = "21NitroExNTX"
$ws.Range("B23").Value = "LEO"
$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$c = $ws.Range("D23")
$c.Value = "'3.540"
$c.ClearFormats()
$ws.Range("E23").Value = "22LEOLEO"
$ws.Range("B24").Value = "BTSEToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$c = $ws.Range("D24")
$c.Value = "'2.158"
$c.ClearFormats()
$ws.Range("E24").Value = "23BTSETokenBTSE"
$c = $ws.Range("D25")
$c.Value = "'0.3200"
$c.ClearFormats()
$c = $ws.Range("D26")
$c.Value = "'0.1337"
$c.ClearFormats()
$c = $ws.Range("D27")
$c.Value = "'0.0002329"
$c.ClearFormats()
$c = $ws.Range("D40")
$c.Value = "'0.03678"
$c.ClearFormats()
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$c = $ws.Range("D41")
$c.Value = "'0.1050"
$c.ClearFormats()
$ws.Range("E41").Value = "40BKEXTokenBKK"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$c = $ws.Range("D42")
$c.Value = "'0.002259"
$c.ClearFormats()
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$c = $ws.Range("D43")
$c.Value = "'0.003055"
$c.ClearFormats()
$ws.Range("E43").Value = "42KickTokenKICK"
$c = $ws.Range("D44")
$c.Value = "'0.008790"
$c.ClearFormats()
$ws.Range("E44").Value = "43LocalTradersLCTBestin24h"
$c = $ws.Range("D45")
$c.Value = "'0.00005485"
$c.ClearFormats()
$c = $ws.Range("D47")
$c.Value = "'0.3999"
$c.ClearFormats()
$c = $ws.Range("D48")
$c.Value = "'0.002560"
$c.ClearFormats()
$c = $ws.Range("D49")
$c.Value = "'0.00002100"
$c.ClearFormats()
$c = $ws.Range("D50")
$c.Value = "'0.0002000"
$c.ClearFormats()
